$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Add Panels")

$ws.Range("F2").Value = "Alarm Current(A)"
$ws.Range("G2").Value = "Standby Current(A)"

$ws.Activate() | Out-Null
$ws.Range("F2:G2").Select() | Out-Null
